# Updated MTOM output files and April model file with run testing new development ruleset
#
# 1) Corrects a typo in the shared header text used in column B1 on every
#    Trace sheet: "DCP BWSCP Flags.LB DCP BWSP" -> "DCP BWSCP Flags.LB DCP BWSCP"
# 2) Updates the re-run model results (columns T/I/O/R on rows 2-3) on each
#    Trace1..Trace38 sheet to the values produced by the new ruleset run.

$wb = $excel.ActiveWorkbook

# Fix typo in shared header string "BWSP" -> "BWSCP" across all Trace sheets (column B1)
for ($i = 1; $i -le 38; $i++) {
    $ws = $wb.Worksheets.Item("Trace$i")
    $ws.Range("B1").Value = "DCP BWSCP Flags.LB DCP BWSCP"
}

# Updated calculation results (re-run with new development ruleset)
$ws = $wb.Worksheets.Item("Trace1")
$ws.Range("T2").Value = 8402103.3535112683

$ws = $wb.Worksheets.Item("Trace10")
$ws.Range("T2").Value = 8448537.9355112668
$ws.Range("I3").Value = 748403.90000073391
$ws.Range("O3").Value = 698403.90000073391

$ws = $wb.Worksheets.Item("Trace11")
$ws.Range("T2").Value = 8425298.5215112679

$ws = $wb.Worksheets.Item("Trace12")
$ws.Range("T2").Value = 8409308.3775112685
$ws.Range("I3").Value = 938522.70000073418
$ws.Range("O3").Value = 888522.70000073418

$ws = $wb.Worksheets.Item("Trace13")
$ws.Range("T2").Value = 8414337.971511269
$ws.Range("I3").Value = 748403.90000073391
$ws.Range("O3").Value = 698403.90000073391

$ws = $wb.Worksheets.Item("Trace14")
$ws.Range("T2").Value = 8435463.8115112688
$ws.Range("I3").Value = 748403.90000073391
$ws.Range("O3").Value = 698403.90000073391

$ws = $wb.Worksheets.Item("Trace15")
$ws.Range("T2").Value = 8460768.0255112685

$ws = $wb.Worksheets.Item("Trace16")
$ws.Range("T2").Value = 8455180.7455112673
$ws.Range("I3").Value = 728264.500000734
$ws.Range("O3").Value = 678264.500000734

$ws = $wb.Worksheets.Item("Trace17")
$ws.Range("T2").Value = 8391746.4835112691

$ws = $wb.Worksheets.Item("Trace18")
$ws.Range("T2").Value = 8458151.709511267
$ws.Range("R3").Value = 10911678.653545434
$ws.Range("T3").Value = 11040197.385545431

$ws = $wb.Worksheets.Item("Trace19")
$ws.Range("T2").Value = 8413963.1435112692
$ws.Range("R3").Value = 11110018.693070415
$ws.Range("T3").Value = 11171690.379070412

$ws = $wb.Worksheets.Item("Trace2")
$ws.Range("T2").Value = 8421280.3535112664
$ws.Range("R3").Value = 11185968.501632852
$ws.Range("T3").Value = 11348426.501632851

$ws = $wb.Worksheets.Item("Trace20")
$ws.Range("T2").Value = 8370714.3735112697

$ws = $wb.Worksheets.Item("Trace21")
$ws.Range("T2").Value = 8325968.4535112707

$ws = $wb.Worksheets.Item("Trace22")
$ws.Range("T2").Value = 8440828.2735112682

$ws = $wb.Worksheets.Item("Trace23")
$ws.Range("T2").Value = 8453200.0635112692

$ws = $wb.Worksheets.Item("Trace24")
$ws.Range("T2").Value = 8397154.1435112692

$ws = $wb.Worksheets.Item("Trace25")
$ws.Range("T2").Value = 8427449.2635112684

$ws = $wb.Worksheets.Item("Trace26")
$ws.Range("T2").Value = 8427691.9135112688

$ws = $wb.Worksheets.Item("Trace27")
$ws.Range("T2").Value = 8377021.3035112685
$ws.Range("R3").Value = 9470964.2728192639
$ws.Range("T3").Value = 9600132.0528192632

$ws = $wb.Worksheets.Item("Trace28")
$ws.Range("T2").Value = 8426573.6935112681

$ws = $wb.Worksheets.Item("Trace29")
$ws.Range("T2").Value = 8487087.0835112706

$ws = $wb.Worksheets.Item("Trace3")
$ws.Range("T2").Value = 8448689.3535112701

$ws = $wb.Worksheets.Item("Trace30")
$ws.Range("T2").Value = 8401880.5535112675
$ws.Range("I3").Value = 748403.90000073391
$ws.Range("O3").Value = 698403.90000073391

$ws = $wb.Worksheets.Item("Trace31")
$ws.Range("T2").Value = 8450398.4335112683

$ws = $wb.Worksheets.Item("Trace32")
$ws.Range("T2").Value = 8425821.3235112689

$ws = $wb.Worksheets.Item("Trace33")
$ws.Range("T2").Value = 8429665.3535112683

$ws = $wb.Worksheets.Item("Trace34")
$ws.Range("T2").Value = 8473276.1335112676

$ws = $wb.Worksheets.Item("Trace35")
$ws.Range("T2").Value = 8368759.9235112704

$ws = $wb.Worksheets.Item("Trace36")
$ws.Range("T2").Value = 8382946.1235112678
$ws.Range("I3").Value = 748403.90000073391
$ws.Range("O3").Value = 698403.90000073391

$ws = $wb.Worksheets.Item("Trace37")
$ws.Range("T2").Value = 8413026.903511269
$ws.Range("I3").Value = 748403.90000073391
$ws.Range("O3").Value = 698403.90000073391

$ws = $wb.Worksheets.Item("Trace38")
$ws.Range("T2").Value = 8415215.3935112692
$ws.Range("I3").Value = 938522.70000073418
$ws.Range("O3").Value = 888522.70000073418

$ws = $wb.Worksheets.Item("Trace4")
$ws.Range("T2").Value = 8413142.1275112685

$ws = $wb.Worksheets.Item("Trace5")
$ws.Range("T2").Value = 8445953.1055112705
$ws.Range("R3").Value = 10237282.93242478
$ws.Range("T3").Value = 10511132.814424779

$ws = $wb.Worksheets.Item("Trace6")
$ws.Range("T2").Value = 8491662.2595112678

$ws = $wb.Worksheets.Item("Trace7")
$ws.Range("T2").Value = 8496967.8915112689
$ws.Range("R3").Value = 10681328.641753327
$ws.Range("T3").Value = 10906231.973753328

$ws = $wb.Worksheets.Item("Trace8")
$ws.Range("T2").Value = 8416779.8535112683

$ws = $wb.Worksheets.Item("Trace9")
$ws.Range("T2").Value = 8381092.5175112681

